$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows right after the header row (becoming rows 2-8),
# shifting the existing data rows (old 2-21) down to rows 9-28.
$ws.Rows("2:8").Insert()
# The inserted rows inherit the header's formatting by default; clear it
# so the new data rows stay unstyled like the rest of the data rows.
$ws.Rows("2:8").ClearFormats()

# Fill the newly inserted rows 2-8 with the new accelerometer samples.
$top = New-Object 'object[,]' 7,3
$top[0,0] = -0.7812347412109375
$top[0,1] = 2.500675392150879
$top[0,2] = 2.085890746116638

$top[1,0] = -0.6001354455947877
$top[1,1] = 2.425701588392258
$top[1,2] = 1.942269176244736

$top[2,0] = -0.4221334457397461
$top[2,1] = 2.401677787303925
$top[2,2] = 1.923809313774109

$top[3,0] = -0.4831114292144776
$top[3,1] = 2.420794081687927
$top[3,2] = 2.053104478120804

$top[4,0] = -0.4990121841430664
$top[4,1] = 2.45754919052124
$top[4,2] = 2.050947427749634

$top[5,0] = -0.4582573175430297
$top[5,1] = 2.450900214910507
$top[5,2] = 1.956571793556213

$top[6,0] = -0.5409791469573973
$top[6,1] = 2.33053719997406
$top[6,2] = 1.996131032705307

$ws.Range("A2:C8").Value2 = $top

# Append 3 new rows of data at the end (rows 29-31).
$bottom = New-Object 'object[,]' 3,3
$bottom[0,0] = 0.1567803621292113
$bottom[0,1] = 2.781254351139069
$bottom[0,2] = 0.9996474064886576

$bottom[1,0] = 0.246018409729004
$bottom[1,1] = 2.646198272705077
$bottom[1,2] = 1.203094172477722

$bottom[2,0] = 0.4407022714614867
$bottom[2,1] = 2.732884711027145
$bottom[2,2] = 1.124532252550125

$ws.Range("A29:C31").Value2 = $bottom
